# Practice tasks and final revisions
# Rename sheets and update task-order CSV filenames.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16509961473895752"
$ws1.Range("B2").Value = "go_stims-16509961473495378.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961473735778.csv"
$ws1.Range("B4").Value = "go_stims-16509961473735778.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961473895752.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1650996149439096"
$ws2.Range("B2").Value = "OB-16509961485348463.csv"
$ws2.Range("B3").Value = "ZB-match_0-16509961473895752.csv"
$ws2.Range("B4").Value = "ZB-match_4-16509961475495794.csv"
$ws2.Range("B5").Value = "TB-16509961494150436.csv"
$ws2.Range("B6").Value = "OB-16509961479908764.csv"
$ws2.Range("B7").Value = "ZB-match_6-1650996147590886.csv"
$ws2.Range("B8").Value = "TB-1650996148991043.csv"
$ws2.Range("B9").Value = "TB-1650996148919078.csv"
$ws2.Range("B10").Value = "OB-16509961481348789.csv"

# --- Sheet 3: RS (name only) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1650996149439096"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509961494870749"
$ws4.Range("B2").Value = "MM_stims-1650996149455077.csv"
$ws4.Range("B3").Value = "ZM_stims-1650996149439096.csv"
$ws4.Range("B4").Value = "MM_stims-16509961494710808.csv"
$ws4.Range("B5").Value = "ZM_stims-1650996149455077.csv"
$ws4.Range("B6").Value = "MM_stims-16509961494870749.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961494710808.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16509961495590568"
$ws5.Range("B2").Value = "vSAT_stims-16509961495350425.csv"
$ws5.Range("B3").Value = "vSAT_stims-165099614951908.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961495030527.csv"
$ws5.Range("B5").Value = "SAT_stims-16509961494870749.csv"
